$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Timeline")
$lo = $ws.ListObjects.Item("Activities")

# Insert a new row before the last one ("Project End"), so the new
# activity becomes Id 3 and "Project End" shifts down to Id 4.
$newRow = $lo.ListRows.Add(13)

$ws.Range("C16").Value = "Abandonment"
$ws.Range("D16").Value = 43801
$ws.Range("E16").Formula = "=Activities[[#This Row],[START]]+2"
$ws.Range("F16").Value = "Ericson's approval to basically do the whole project differently and alone"

$excel.Calculate()
